$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new weekly rows above the current row 556, pushing the existing
# rows 556:567 down to 562:573 (dimension grows from A1:R567 to A1:R573).
$ws.Range("A556:R561").EntireRow.Insert()

# Populate the 6 newly inserted rows (556-561) with the new week's data
# (Fecha = 44448). Columns A,B,C,E,F,G,H,N,Q,R are constant across every
# row of this sheet.

# Row 556
$ws.Range("A556").Value = 9
$ws.Range("B556").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C556").Value = "Metropolitana"
$ws.Range("D556").Value = 44448
$ws.Range("E556").Value = 13
$ws.Range("F556").Value = 100114013
$ws.Range("G556").Value = "Zanahoria"
$ws.Range("H556").Value = "Sin especificar"
$ws.Range("I556").Value = "Primera"
$ws.Range("J556").Value = 250
$ws.Range("K556").Value = 5000
$ws.Range("L556").Value = 5500
$ws.Range("M556").Value = 5250
$ws.Range("N556").Value = "$/saco 20 kilos"
$ws.Range("O556").Value = "Chillán"
$ws.Range("P556").Value = 262
$ws.Range("Q556").Value = 20
$ws.Range("R556").Value = "Hortaliza"

# Row 557
$ws.Range("A557").Value = 9
$ws.Range("B557").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C557").Value = "Metropolitana"
$ws.Range("D557").Value = 44448
$ws.Range("E557").Value = 13
$ws.Range("F557").Value = 100114013
$ws.Range("G557").Value = "Zanahoria"
$ws.Range("H557").Value = "Sin especificar"
$ws.Range("I557").Value = "Primera"
$ws.Range("J557").Value = 160
$ws.Range("K557").Value = 5500
$ws.Range("L557").Value = 6000
$ws.Range("M557").Value = 5750
$ws.Range("N557").Value = "$/saco 20 kilos"
$ws.Range("O557").Value = "Región Metropolitana"
$ws.Range("P557").Value = 288
$ws.Range("Q557").Value = 20
$ws.Range("R557").Value = "Hortaliza"

# Row 558
$ws.Range("A558").Value = 9
$ws.Range("B558").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C558").Value = "Metropolitana"
$ws.Range("D558").Value = 44448
$ws.Range("E558").Value = 13
$ws.Range("F558").Value = 100114013
$ws.Range("G558").Value = "Zanahoria"
$ws.Range("H558").Value = "Sin especificar"
$ws.Range("I558").Value = "Primera"
$ws.Range("J558").Value = 210
$ws.Range("K558").Value = 5000
$ws.Range("L558").Value = 5500
$ws.Range("M558").Value = 5250
$ws.Range("N558").Value = "$/saco 20 kilos"
$ws.Range("O558").Value = "Región de La Araucanía"
$ws.Range("P558").Value = 262
$ws.Range("Q558").Value = 20
$ws.Range("R558").Value = "Hortaliza"

# Row 559
$ws.Range("A559").Value = 9
$ws.Range("B559").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C559").Value = "Metropolitana"
$ws.Range("D559").Value = 44448
$ws.Range("E559").Value = 13
$ws.Range("F559").Value = 100114013
$ws.Range("G559").Value = "Zanahoria"
$ws.Range("H559").Value = "Sin especificar"
$ws.Range("I559").Value = "Segunda"
$ws.Range("J559").Value = 160
$ws.Range("K559").Value = 4000
$ws.Range("L559").Value = 4500
$ws.Range("M559").Value = 4250
$ws.Range("N559").Value = "$/saco 20 kilos"
$ws.Range("O559").Value = "Chillán"
$ws.Range("P559").Value = 212
$ws.Range("Q559").Value = 20
$ws.Range("R559").Value = "Hortaliza"

# Row 560
$ws.Range("A560").Value = 9
$ws.Range("B560").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C560").Value = "Metropolitana"
$ws.Range("D560").Value = 44448
$ws.Range("E560").Value = 13
$ws.Range("F560").Value = 100114013
$ws.Range("G560").Value = "Zanahoria"
$ws.Range("H560").Value = "Sin especificar"
$ws.Range("I560").Value = "Segunda"
$ws.Range("J560").Value = 52
$ws.Range("K560").Value = 4500
$ws.Range("L560").Value = 5000
$ws.Range("M560").Value = 4750
$ws.Range("N560").Value = "$/saco 20 kilos"
$ws.Range("O560").Value = "Región Metropolitana"
$ws.Range("P560").Value = 238
$ws.Range("Q560").Value = 20
$ws.Range("R560").Value = "Hortaliza"

# Row 561
$ws.Range("A561").Value = 9
$ws.Range("B561").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C561").Value = "Metropolitana"
$ws.Range("D561").Value = 44448
$ws.Range("E561").Value = 13
$ws.Range("F561").Value = 100114013
$ws.Range("G561").Value = "Zanahoria"
$ws.Range("H561").Value = "Sin especificar"
$ws.Range("I561").Value = "Segunda"
$ws.Range("J561").Value = 106
$ws.Range("K561").Value = 4000
$ws.Range("L561").Value = 4500
$ws.Range("M561").Value = 4250
$ws.Range("N561").Value = "$/saco 20 kilos"
$ws.Range("O561").Value = "Región de La Araucanía"
$ws.Range("P561").Value = 212
$ws.Range("Q561").Value = 20
$ws.Range("R561").Value = "Hortaliza"
